$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-04 -> 2023-10-05, serial 45203 -> 45204) for every data row
# (rows 2 through 421).
$ws.Range("C2:C421").Value = 45204
